$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "ubicacion" sheet: add a new row describing the source file for
#    the run-off results, right above the first-round row.
# ------------------------------------------------------------------
$wsUbic = $wb.Worksheets.Item("ubicacion")
$wsUbic.Rows("2:2").Copy()
$wsUbic.Rows("2:2").Insert()
$wsUbic.Range("B2").Value = "Resultados_2021_Mesa_PRESIDENCIAL_Tricel_2v_TEMP.xlsx"
$wsUbic.Range("C2").Value = "2021_presidencial_2v "

# ------------------------------------------------------------------
# 2) "meta" sheet: add a new column for the 2021 presidential run-off
#    (2021_presidencial_2v). Duplicate the existing 2021_presidencial_1v
#    column (currently column C) into a new column C, pushing the old
#    column C (and everything after it) one column to the right, then
#    retitle the new column's header.
# ------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("meta")
$wsMeta.Columns("C:C").Copy()
$wsMeta.Columns("C:C").Insert()
$wsMeta.Range("C1").Value = "2021_presidencial_2v "
$wsMeta.Range("E14").Select()

$wsUbic.Range("A13").Select()

# ------------------------------------------------------------------
# 3) "tendencia" sheet: add the run-off candidate options (second
#    round only features the two top candidates from the first round)
# ------------------------------------------------------------------
$wsTend = $wb.Worksheets.Item("tendencia")
$wsTend.Range("A44").Value = "2021_presidencial_2v "
$wsTend.Range("B44").Value = "opcion"
$wsTend.Range("C44").Value = "1. GABRIEL BORIC FONT"
$wsTend.Range("D44").Value = -1

$wsTend.Range("A45").Value = "2021_presidencial_2v "
$wsTend.Range("B45").Value = "opcion"
$wsTend.Range("C45").Value = "2. JOSE ANTONIO KAST RIST"
$wsTend.Range("D45").Value = 1

$wsTend.Range("A46").Value = "2021_presidencial_2v "
$wsTend.Range("B46").Value = "opcion"
$wsTend.Range("C46").Value = "Votos Nulos"

$wsTend.Range("A47").Value = "2021_presidencial_2v "
$wsTend.Range("B47").Value = "opcion"
$wsTend.Range("C47").Value = "Votos en Blanco"

$wsTend.Range("D46").Select()
